# Add 2022-Q4 data
# 1) Insert a new worksheet "2022-Q4" right before the existing "2022-Q2" sheet,
#    populate it with the fund-holding data for the new quarter.
# 2) Update the "总计" (summary) sheet: insert a new row for 2022-Q4 at the top
#    of the data (row 2), shifting the existing quarters down by one row, and
#    renumber the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet, positioned before "2022-Q2".
# Duplicate the existing "2022-Q2" sheet (so sheet-level formatting such as
# sheetPr/pageMargins/header style is preserved exactly), trim it down to a
# single data row, then overwrite the values with the new quarter's data.
# ---------------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$newSheet = $wb.Worksheets.Item("2022-Q2 (2)")
$newSheet.Name = "2022-Q4"

# The duplicated sheet has 3 data rows (rows 2-4); only one is needed.
$newSheet.Rows.Item(3).Resize(2).Delete()

# Data row 2 (single holding for this quarter) - header row is already correct
$newSheet.Range("A2").Value = 0

$newSheet.Range("B2").Value = "'159804"
$newSheet.Range("B2").ClearFormats()

$newSheet.Range("C2").Value = "国寿安保国证创业板中盘精选88ETF"

$newSheet.Range("D2").Value = "'1.15"
$newSheet.Range("D2").ClearFormats()

$newSheet.Range("E2").Value = "'99.00"
$newSheet.Range("E2").ClearFormats()

$newSheet.Range("F2").Value = "'1.96"
$newSheet.Range("F2").ClearFormats()

$newSheet.Range("G2").Value = "'0.0225"
$newSheet.Range("G2").ClearFormats()

$newSheet.Range("H2").Value = 6

# ---------------------------------------------------------------------------
# Step 2: update the "总计" summary sheet
# ---------------------------------------------------------------------------
$tot = $wb.Worksheets.Item("总计")

# Give the new row-2 index cell (A2) the same look as the existing index cells
$tot.Range("A8").Copy()
$tot.Range("A9").PasteSpecial(-4122)

# Shift the existing quarter rows (2..8) down to (3..9), processing from the
# bottom up so we never overwrite data we still need to read.
for ($r = 8; $r -ge 2; $r--) {
    $dst = $r + 1
    $tot.Cells.Item($dst, 2).Value = $tot.Cells.Item($r, 2).Value2
    $tot.Cells.Item($dst, 3).Value = $tot.Cells.Item($r, 3).Value2
    $tot.Cells.Item($dst, 4).Value = $tot.Cells.Item($r, 4).Value2
}

# New row 2: 2022-Q4 summary data
$tot.Cells.Item(2, 2).Value = "2022-Q4"
$tot.Cells.Item(2, 3).Value = 1
$tot.Cells.Item(2, 4).Value = 0.02

# Renumber the index column (A2:A9 -> 0..7)
for ($r = 2; $r -le 9; $r++) {
    $tot.Cells.Item($r, 1).Value = $r - 2
}

# Restore "总计" as the active sheet (matches the workbook's original state)
$tot.Activate()
$tot.Range("A1").Select() | Out-Null

Write-Output "2022-Q4 sheet added and 总计 summary updated."
